$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (2023-09-24)
$ws.Range("B2").Value = 0.06328177979961902
$ws.Range("C2").Value = 0.05231270169004087
$ws.Range("D2").Value = 3.082599426703578
$ws.Range("E2").Value = 246.9852506941017
$ws.Range("G2").Value = 250.1834446022949

# Row 3 (2023-03-09)
$ws.Range("B3").Value = 0.1554434735375247
$ws.Range("C3").Value = 0.05231270169004087
$ws.Range("D3").Value = 3.082599426703578
$ws.Range("E3").Value = 6.48142807727062
$ws.Range("G3").Value = 9.771783679201764
